$d = $word.ActiveDocument
$paras = $d.Paragraphs

# ---------------------------------------------------------------------------
# Locate the two target paragraphs by content instead of a hard-coded index,
# so the script is resilient to any incidental paragraph-numbering drift.
# ---------------------------------------------------------------------------
$bulletPara = $null
$demoPara = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    $t = $p.Range.Text
    if ($t.TrimEnd("`r", "`n") -eq " ") {
        # First lone-space placeholder paragraph in the doc is the one right
        # under "Distribution Channels:".
        if ($null -eq $bulletPara) { $bulletPara = $p }
    }
    if ($t.StartsWith("For our Demographics analysis")) {
        $demoPara = $p
    }
}

# ---------------------------------------------------------------------------
# Change 1: "Distribution Channels" bullet — replace the placeholder single
# space with the bar-plot design explanation paragraph.
# ---------------------------------------------------------------------------
$barPlotText = "This bar plot was intended to uncover the variation in cancellation rates through booking channels without overly simplifying the visualization for viewers. We used a sorted bar order (top-to-bottom by highest to lowest cancellation rate) to prioritize attention towards the most at risk for volatility. First, all the bars are rendered as identical blue colors for consistency of look. A rollover hover interaction tailored in orange color helps the user differentiate and identify specific segments. We used clear labeling of axes, chart title, and percentage scaling on the Y-axis for ease of interpretation. Overall, the design strikes a balance between interactivity and legibility such that users can contrast channel cancellation rates simply immediately and easily understand the patterns."
$bulletRange = $bulletPara.Range
$bulletRange.Find.Execute(" ", $true, $false, $false, $false, $false, $true, 1, $false, $barPlotText, 2)

# ---------------------------------------------------------------------------
# Change 2: "Demographics" bullet — drop the grammar-check (proofErr) run
# splits, merging the text back into two runs, with a rendered-page-break
# marker inserted right before "them to assess...".
# ---------------------------------------------------------------------------
$demoFull = $demoPara.Range
# Exclude the trailing paragraph mark so InsertXML only rewrites the runs,
# leaving the paragraph's pPr (ListParagraph style/numbering) untouched.
$demoRange = $d.Range($demoFull.Start, $demoFull.End - 1)

$run1 = @'
For our Demographics analysis, we wanted to hone in on how price might influences per the specific type of traveler (solo, family, or couple). The way we felt would be a best representation of this was if we utilized the Average Daily Room Rate (ADR) &#8211; which is representative of the price that guests have to pay &#8211; and to aggregate it to a monthly level. After this aggregation, we wanted to plot it over all monthly data in the dataset to be able to see apparent trends. Next, another visualization would be added to show how different traveler group types cancelled per the specific time period; we wanted to hone in on each group&#8217;s percentage of the total amount, to indicate prevalency/significance. After this, we wanted to include an interactive element in this visualization by allowing users to select filters, perhaps during price spikes or troughs, which will allow 
'@

$run2 = @'
them to assess cancelation rates (per traveler type) for different time periods. This is an effective feature of this visualization, as it allows the user to potentially see nuances in the data if they select a smaller time period; however, it&#8217;s beneficial to our conclusion because, after doing so, one could see that couples consistently have the highest cancellation rates.
'@

$run1 = $run1.TrimEnd("`r", "`n")
$run2 = $run2.TrimEnd("`r", "`n")

$xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p>' +
  '<w:r><w:t xml:space="preserve">' + $run1 + '</w:t></w:r>' +
  '<w:r><w:lastRenderedPageBreak/><w:t>' + $run2 + '</w:t></w:r>' +
  '</w:p>' +
  '</w:body>' +
  '</w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$demoRange.InsertXML($xmlFrag)

Write-Output "done"
